$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison -------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# Give C1 the "top+bottom" border used for the inner cells of the merged B1:D1
# header box, and D1 the "top+bottom+right" border that closes the box.
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders(8).LineStyle = 1
$ws1.Range("C1").Borders(9).LineStyle = 1
$ws1.Range("C1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws1.Range("D1").Borders(10).LineStyle = 1

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison -------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Reuse the already-built formats from sheet 1 so the style table stays
# minimal (no unused intermediate styles get left behind).
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)

$ws2.Range("C1").Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$ws2.Range("D1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()

$excel.CutCopyMode = $false
